# netCrypto.xlsx — "Add files via upload" re-save.
#
# Observed changes (vs. before.xlsx):
#   - Worksheet scroll position: topLeftCell L1 -> H1
#   - Worksheet selection:       R13 -> T3
#   - Cell T2 value:              70297 -> 85402
#
# (File metadata such as fileVersion/rupBuild, the x15ac absPath hint, and
# the xr:revisionPtr coauthoring GUID are Excel/host bookkeeping written by
# the save pipeline itself, not values exposed on the Excel object model,
# so there is nothing for workbook-editing script to set for those.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll the view so column H is the left-most visible column, row 1 on top.
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1

# Move the selection/active cell to T3.
$ws.Range("T3").Select()

# Update the USD Amount value in T2.
$ws.Range("T2").Value = 85402
